# ---------------------------------------------------------------------------
# "last minute changes before sending to arno"
#
# 1. Minor cosmetic tidy-up of Sheet1 (selection, a couple of re-sized
#    columns).
# 2. A brand-new worksheet "Actuations when Even" is added right after
#    Sheet1 and becomes the active tab. It re-derives the same embodied /
#    operational-carbon break-even table as Sheet1, but only keeps the two
#    scenarios that matter here (columns B = soft robotic actuators,
#    F = motorized / servo actuators), adds a title, an "Input" cell for the
#    number of actuations, and a new "Difference" output row comparing the
#    two scenarios' gCO2/kWh.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet1: small formatting touch-ups
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Columns.Item(5).ColumnWidth = 17
$ws1.Columns.Item(6).ColumnWidth = 30.6640625
$ws1.Columns.Item(7).ColumnWidth = 25.83203125

$ws1.Range("A1:XFD1048576").Select()

# ---------------------------------------------------------------------------
# Add the new sheet right after Sheet1
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Actuations when Even"

$ws2.Columns.Item(2).ColumnWidth = 22.83203125
$ws2.Columns.Item(3).ColumnWidth = 23
$ws2.Columns.Item(4).ColumnWidth = 18.1640625
$ws2.Columns.Item(5).ColumnWidth = 17
$ws2.Columns.Item(6).ColumnWidth = 30.6640625
$ws2.Columns.Item(7).ColumnWidth = 25.83203125
$ws2.Columns.Item(8).ColumnWidth = 32.33203125

# NOTE: shared-string table entries are assigned in first-seen order, so we
# touch the three brand-new text labels in the same order the original
# author must have (Number of Actuations -> title -> Difference) before
# filling in the rest of the sheet.
$ws2.Range("A18").Value = "Number of  Actuations"
$ws2.Range("A1").Value = "Find out when soft robotic actuators break even with servo motors"
$ws2.Range("A19").Value = "Difference"

# ---- Row 1: title ---------------------------------------------------------
$ws2.Range("A1").Font.Bold = $true

# ---- Row 2: headers --------------------------------------------------------
$ws2.Range("B2").Value = "ASF ENTSO-E incl shading"
$ws2.Range("F2").Value = "ASF ENTSO-E incl shading motorized"

# ---- Row 3: Embodied -------------------------------------------------------
$ws2.Range("A3").Value = "Embodied"
$ws2.Range("B3").Value = 2675.4
$ws2.Range("F3").Value = 3251.2
$ws2.Range("B3:H3").NumberFormat = "0.0"

# ---- Row 4: HVAC offset ----------------------------------------------------
$ws2.Range("A4").Value = "HVAC offset"
$ws2.Range("B4").Formula = "=-1021.8*20*0.4621"
$ws2.Range("F4").Formula = "=-1021.8*20*0.4621"
$ws2.Range("B4:H4").NumberFormat = "0.0"

# ---- Row 5: Actuators -------------------------------------------------------
$ws2.Range("A5").Value = "Actuators"
$ws2.Range("B5").Formula = "=0.31/1000*`$B`$18*54*365*20*0.4621"
$ws2.Range("F5").Formula = "=0.01/1000*`$B`$18*54*365*20*0.4621"
$ws2.Range("B5:H5").NumberFormat = "0.0"

# ---- Row 6: Maintenance -----------------------------------------------------
$ws2.Range("A6").Value = "Maintenance"
$ws2.Range("B6").Formula = "=(42.73+35.81)*3"
$ws2.Range("F6").Value = 0
$ws2.Range("B6:H6").NumberFormat = "0.0"

# ---- Row 7: Disposal --------------------------------------------------------
$ws2.Range("A7").Value = "Disposal"
$ws2.Range("B7").Value = 77
$ws2.Range("F7").Value = 86.5
$ws2.Range("B7:H7").NumberFormat = "0.0"

# ---- Row 8: Total (bold) -----------------------------------------------------
$ws2.Range("A8").Value = "Total"
$ws2.Range("A8").Font.Bold = $true
$ws2.Range("B8").Formula = "=SUM(B3:B7)"
$ws2.Range("F8").Formula = "=SUM(F3:F7)"
$ws2.Range("B8:H8").NumberFormat = "0.0"
$ws2.Range("B8:H8").Font.Bold = $true

# ---- Row 9: spacer (bold, empty) ---------------------------------------------
$ws2.Range("A9").Font.Bold = $true

# ---- Row 10: Elec prod. (italic) ---------------------------------------------
$ws2.Range("A10").Value = "Elec prod."
$ws2.Range("A10").Font.Italic = $true
$ws2.Range("B10").Formula = "=580*20"
$ws2.Range("F10").Formula = "=580*20"
$ws2.Range("B10:H10").Font.Italic = $true

# ---- Row 12: gCO2/kWh (bold) ---------------------------------------------------
$ws2.Range("A12").Value = "gCO2/kWh"
$ws2.Range("A12").Font.Bold = $true
$ws2.Range("B12").Formula = "=B8/B10*1000"
$ws2.Range("F12").Formula = "=F8/F10*1000"
$ws2.Range("B12:H12").NumberFormat = "0.0"
$ws2.Range("B12:H12").Font.Bold = $true

# ---- Row 16: calc for numb of actuations ----------------------------------------
$ws2.Range("A16").Value = "calc for numb of actuations"
$ws2.Range("B16").Formula = "=0.31/1000*`$B`$18*54*365*20*0.4621"
$ws2.Range("F16").Formula = "=F8-B8"
$ws2.Range("B16:D16").NumberFormat = "0.0"
$ws2.Range("F16").NumberFormat = "0.0"
$ws2.Range("H16").NumberFormat = "0.0"

# ---- Row 18: Number of Actuations (input) ---------------------------------------
$ws2.Range("A18").Font.Bold = $true
$ws2.Range("B18").Value = 6
$ws2.Range("B18").Style = "Input"

# ---- Row 19: Difference (output) -------------------------------------------------
$ws2.Range("A19").Font.Bold = $true
$ws2.Range("B19").Formula = "=F12-B12"
$ws2.Range("B19").Style = "Output"
$ws2.Range("B19").NumberFormat = "0.0"

$ws2.Range("H26").Select()
